# Update odds values on Sheet1 to reflect the latest Betfair Back/Lay snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Melbourne City vs Brisbane Roar)
$ws.Range("F2").Value = 1.86
$ws.Range("G2").Value = 1.87
$ws.Range("H2").Value = 4.9
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 3.8
$ws.Range("L2").Value = 1.45
$ws.Range("N2").Value = 3.35
$ws.Range("P2").Value = 1.77
$ws.Range("Q2").Value = 2.16
$ws.Range("V2").Value = 1.25
$ws.Range("W2").Value = 2.14
$ws.Range("AB2").Value = 8
$ws.Range("AG2").Value = 10.5
$ws.Range("AO2").Value = 90

# Row 3 (Pisa vs Como)
$ws.Range("H3").Value = 1.8
$ws.Range("I3").Value = 1.82
$ws.Range("J3").Value = 3.65
$ws.Range("K3").Value = 3.7
$ws.Range("N3").Value = 3.35
$ws.Range("P3").Value = 1.78
$ws.Range("Z3").Value = 9.6
$ws.Range("AB3").Value = 16.5
$ws.Range("AO3").Value = 14.5

# Row 4 (Lecce vs Roma)
$ws.Range("J4").Value = 3.55
$ws.Range("L4").Value = 1.52

# Row 5 (Sassuolo vs Juventus)
$ws.Range("F5").Value = 5.3
$ws.Range("J5").Value = 3.9
$ws.Range("K5").Value = 3.95
$ws.Range("P5").Value = 1.9
$ws.Range("R5").Value = 1.35

# Row 6 (Livingston vs St Mirren)
$ws.Range("F6").Value = 2.66
$ws.Range("G6").Value = 2.7
$ws.Range("H6").Value = 2.76
$ws.Range("I6").Value = 2.78
$ws.Range("J6").Value = 3.65
$ws.Range("K6").Value = 3.75
$ws.Range("L6").Value = 1.4
$ws.Range("N6").Value = 3.45
$ws.Range("P6").Value = 1.85
$ws.Range("Q6").Value = 1.96
$ws.Range("T6").Value = 1.79
$ws.Range("U6").Value = 2.04
$ws.Range("V6").Value = 1.56
$ws.Range("W6").Value = 1.57
$ws.Range("X6").Value = 17
$ws.Range("Y6").Value = 14
$ws.Range("AA6").Value = 50
$ws.Range("AB6").Value = 11
$ws.Range("AC6").Value = 9
$ws.Range("AD6").Value = 13
$ws.Range("AF6").Value = 18
$ws.Range("AG6").Value = 12
$ws.Range("AH6").Value = 20
$ws.Range("AJ6").Value = 44
$ws.Range("AK6").Value = 980
$ws.Range("AM6").Value = 120
$ws.Range("AN6").Value = 28
$ws.Range("AO6").Value = 27

# Row 7 (West Ham vs Nottm Forest)
$ws.Range("F7").Value = 3.3
$ws.Range("G7").Value = 3.35
$ws.Range("AJ7").Value = 60

# Row 8 (Rangers vs Aberdeen)
$ws.Range("G8").Value = 1.66
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 6.6
$ws.Range("J8").Value = 4.2
$ws.Range("K8").Value = 4.5
$ws.Range("N8").Value = 4.5
$ws.Range("O8").Value = 1.24
$ws.Range("P8").Value = 2.2
$ws.Range("R8").Value = 1.46
$ws.Range("S8").Value = 2.64
$ws.Range("T8").Value = 1.04
$ws.Range("U8").Value = 2.06
$ws.Range("V8").Value = 1.17
$ws.Range("W8").Value = 2.5
$ws.Range("X8").Value = 1000
$ws.Range("Z8").Value = 1000
$ws.Range("AA8").Value = 1000
$ws.Range("AB8").Value = 10
$ws.Range("AC8").Value = 1000
$ws.Range("AD8").Value = 24
$ws.Range("AE8").Value = 1000
$ws.Range("AG8").Value = 9.6
$ws.Range("AH8").Value = 20
$ws.Range("AL8").Value = 1000
$ws.Range("AM8").Value = 1000
$ws.Range("AN8").Value = 8.6
$ws.Range("AO8").Value = 1000
